$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 19, whose phone number (column A) is
# stored as text "09876543" (leading zero). The edit appends a brand new
# row 20 that keeps that original text value verbatim (same blank
# birthday, same 0 points), while row 19 itself gets its phone value
# turned into a real number (9876543 - the leading zero is dropped).
#
# Insert the new row *after* row 19 (i.e. at row 20) rather than shifting
# row 19 down - that way row 19's untouched cells (its blank birthday /
# its 0 points) are never rewritten, so they keep their original storage
# exactly as-is, and only column A of row 19 is touched.
$ws.Rows("20:20").Insert()

# New row 20 = verbatim copy of the old row 19 content.
# Force column A to stay text (it looks numeric, so it would otherwise be
# auto-converted) by setting a text number format before assigning it,
# then drop back to the workbook's default cell style so no stray
# formatting is left behind on the cell.
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "09876543"
$ws.Range("A20").Style = "Normal"
$ws.Range("C20").Value = 0

# Row 19's phone becomes a plain number; its other cells are left alone.
$ws.Range("A19").Value = 9876543
